# Rewrite 0.04: II_Windows UI - Recreated UI for Patient Editor
# Adds a new "Devices" section (Devices/Ventilator/IABP/Defibrillator/
# 12LeadECG/IVPump/LabResults/Cardiotocograph) above the existing
# "DeviceOptions..EditPatient" block on Sheet1, pushing that block down
# by ten rows (33-43 -> 43-53) and inserting the new block at rows 31-39
# (CardiacMonitor, which used to be the lone row 31 entry, becomes row 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: push the existing "DeviceOptions .. EditPatient" block
#     (old rows 33-43) down to rows 43-53. Write from the bottom up so a
#     literal-value write never depends on a cell that still needs moving.
$ws.Range("A53").Value = "EditPatient"
$ws.Range("A52").Value = "NewPatient"
$ws.Range("A51").Value = "PatientOptions"
$ws.Range("A50").Value = "CloseDevice"
$ws.Range("A49").Value = "ToggleFullscreen"
$ws.Range("A48").Value = "ColorScheme"
$ws.Range("A47").Value = "FontSize"
$ws.Range("A46").Value = "TracingRowAmounts"
$ws.Range("A45").Value = "NumericRowAmounts"
$ws.Range("A44").Value = "PauseDevice"
$ws.Range("A43").Value = "DeviceOptions"

# --- Step 2: the old lone row 31 ("CardiacMonitor") moves down to row 32.
$ws.Range("A32").Value = "CardiacMonitor"

# --- Step 3: write the new "Devices" section into rows 31, 33-39.
#     Order matters here: each brand-new string gets interned into the
#     shared-string table in first-write order, and the target file needs
#     them appended as Devices, Ventilator, IABP, Defibrillator, 12LeadECG,
#     IVPump, LabResults, Cardiotocograph (in that order) - independent of
#     which row number each one ends up on.
$ws.Range("A31").Value = "Devices"
$ws.Range("A35").Value = "Ventilator"
$ws.Range("A36").Value = "IABP"
$ws.Range("A34").Value = "Defibrillator"
$ws.Range("A33").Value = "12LeadECG"
$ws.Range("A38").Value = "IVPump"
$ws.Range("A39").Value = "LabResults"
$ws.Range("A37").Value = "Cardiotocograph"

# --- Step 4: rows 40-42 held the tail of the old block that has now been
#     relocated below row 43 - clear them so they go back to being blank.
$ws.Range("A40:A42").ClearContents()

# --- Step 5: move the active selection to A37, matching the new layout.
$ws.Range("A37").Select() | Out-Null
